$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before Q for the 'dload' values; this shifts
# existing columns Q:V (w, piezo_y, hw, u, phi, c) to R:W.
$ws.Columns("Q").Insert()

# New column header
$ws.Cells.Item(1, 17).Value = "dload"

$data = @(
    @(1, 60.52735684831512, 83.9998059745921, 84, 70.39551763623633, 61.19606926906113, 84, 65.46143724227572, 70.99690551636458, 84, 9.86816078792122, 65.89820373474561, 24.1654007820007, 13.00309448363542, 0, 0, 0, 16681.16152365595, 79.12718083676965, 8.130275320405076, 507.3291799932767, 28, 200),
    @(2, 70.39551763623633, 61.19606926906113, 84, 80.26367842415756, 46.34520128139259, 84, 75.32959803019695, 53.17327651921067, 84, 9.868160787921227, 56.15905274156578, 17.72014640795434, 20, 10.82672348078933, 0, 0, 38477.99982237014, 78.97327583538943, 25.79999931617876, 1609.919957329554, 28, 200),
    @(3, 80.26367842415756, 46.34520128139259, 84, 90.13183921207877, 35.16253607942096, 84, 85.19775881811816, 40.40378477144077, 84, 9.868160787921212, 48.44782535860768, 14.87733183099115, 20, 23.59621522855923, 0, 0, 53599.36755201326, 78.17315469042285, 37.76936991898208, 2356.808682944481, 28, 200),
    @(4, 90.13183921207877, 35.16253607942096, 84, 100, 26.32503002402403, 84, 95.06591960603939, 30.4983643078416, 84, 9.868160787921227, 41.768101924379, 13.23081183140688, 20, 24, 9.501635692158402, 0, 66454.32539392641, 77.37303354545627, 46.87466923761467, 2924.979360427155, 28, 200),
    @(5, 100, 26.32503002402403, 84, 110, 19.12879499083992, 84, 105, 22.53205655191036, 84, 10, 35.68533492922602, 12.31174025844438, 20, 24, 17.46794344808964, 25000, 102857.6853514783, 76.56756756756756, 54.03551101565721, 3371.815887377009, 28, 200),
    @(6, 110, 19.12879499083992, 84, 120, 13.3463549614923, 84, 115, 16.07695154586736, 84, 10, 30.00000015313622, 11.54700540161072, 20, 24, 23.92304845413264, 25000, 111378.4239594551, 75.35714285714286, 59.2801913112755, 3699.083937823591, 28, 200),
    @(7, 120, 13.3463549614923, 84, 130, 8.757022693565062, 84, 125, 10.91287885364285, 84, 10, 24.62431846252714, 11.00038197405046, 20, 24, 29.08712114635715, 25000, 118194.9999131914, 73.21428571428571, 62.30140686064286, 3887.607788104114, 28, 200),
    @(8, 130, 8.757022693565062, 84, 140, 5.217597167510036, 84, 135, 6.862915010152392, 84, 10, 19.47122071361051, 10.60660172297659, 20, 24, 33.13708498984761, 25000, 123540.9521865988, 71.07142857142857, 64.20851356127618, 4006.611246223633, 28, 200),
    @(9, 140, 5.217597167510036, 84, 150, 2.633054056944971, 84, 145, 3.810499613777495, 84, 10, 14.47751224073088, 10.327955592437, 20, 24, 36.18950038622251, 25000, 127570.1405098137, 68.5, 64.6895003862225, 4036.624824100284, 28, 200),
    @(10, 150, 2.633054056944971, 84, 162.35, 0.6686231538410539, 74, 156.175, 1.485784080558503, 78.99999999999999, 12.34999999999999, 9.025556359556186, 12.50482896617596, 14.99999999999999, 24, 38.5142159194415, 0, 122436.3747918735, 65.1475, 63.66171591944149, 3972.491073373149, 28, 200),
    @(11, 162.35, 0.6686231538410539, 74, 174.7, 0.000375000585933094, 64, 168.525, 0.17481744224213, 69.00000000000001, 12.34999999999999, 3.093086941691099, 12.3680178772166, 5.000000000000014, 24, 39.82518255775787, 0, 108518.5126056569, 61.44250000000001, 61.26768255775788, 3823.103391604091, 28, 200),
    @(12, 174.7, 0.000375000585933094, 64, 184.5666666666667, 0.3819458071278916, 56, 179.6333333333333, 0.08948243701796343, 60, 9.866666666666674, -2.212803763689251, 9.874029601984471, 0, 20, 39.91051756298204, 0, 75659.45807402786, 56.55384615384616, 56.4643637168282, 3523.376295930079, 32, 100),
    @(13, 184.5666666666667, 0.3819458071278916, 56, 194.4333333333333, 1.58401478028587, 48, 189.5, 0.8792629304200403, 52, 9.866666666666674, -6.940199065150054, 9.939495248171937, 0, 12, 39.12073706957996, 0, 65158.84795942099, 52, 51.12073706957996, 3189.933993141789, 32, 100),
    @(14, 194.4333333333333, 1.58401478028587, 48, 204.3, 3.632006118520721, 40, 199.3666666666667, 2.499933806165217, 44, 9.866666666666674, -11.71569679077994, 10.07659006958593, 0, 4, 37.50006619383478, 0, 53576.08621085047, 44, 41.50006619383478, 2589.604130495291, 32, 100),
    @(15, 204.3, 3.632006118520721, 40, 214.323780536418, 6.626104043639586, 40, 209.311890268209, 5.010025714315489, 40, 10.02378053641795, -16.61459186508567, 10.46050903348247, 0, 0, 34.98997428568451, 0, 46296.60066432835, 40, 34.98997428568451, 2183.374395426713, 27, 400),
    @(16, 214.323780536418, 6.626104043639586, 40, 224.3475610728359, 10.61618850962118, 40, 219.3356708046269, 8.49059100549519, 40, 10.02378053641797, -21.68255994558514, 10.7870149791044, 0, 0, 31.50940899450481, 0, 41691.3288782959, 40, 31.50940899450481, 1966.1871212571, 27, 400),
    @(17, 224.3475610728359, 10.61618850962118, 40, 234.3713416092539, 15.71652194370731, 40, 229.3594513410449, 13.0184593030154, 40, 10.02378053641795, -26.93603065199318, 11.24356274821377, 0, 0, 26.9815406969846, 0, 35700.32960749238, 40, 26.9815406969846, 1683.648139491839, 27, 400),
    @(18, 234.3713416092539, 15.71652194370731, 40, 244.3951221456718, 22.10047486127786, 40, 239.3832318774628, 18.73401630847181, 40, 10.02378053641795, -32.44752375024643, 11.87816108521438, 0, 0, 21.26598369152819, 0, 28137.85305076967, 40, 21.26598369152819, 1326.997382351359, 27, 400),
    @(19, 244.3951221456718, 22.10047486127786, 40, 254.4189026820898, 30.04090987136011, 40, 249.4070124138808, 25.8533245215712, 40, 10.02378053641795, -38.32040221804986, 12.77637964905471, 0, 0, 14.1466754784288, 0, 18718.01848167183, 40, 14.1466754784288, 882.7525498539574, 27, 400),
    @(20, 254.4189026820898, 30.04090987136011, 40, 264.4426832185077, 39.9999598832996, 40, 259.4307929502987, 34.72725405040724, 40, 10.02378053641797, -44.71574021543523, 14.1059568222393, 0, 0, 5.27274594959276, 0, 6976.57596543661, 40, 5.27274594959276, 329.0193472545882, 27, 400)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

Write-Host "Done. Dimension now: $($ws.UsedRange.Address())"
